$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Zahra", "unknown"),
    @("Agus", "Nori"),
    @("Prasetya", "unknown"),
    @("Abil", "unknown"),
    @("Dinul", "unknown"),
    @("Arinal", "unknown"),
    @("Said", "Beling"),
    @("Fajar", "unknown"),
    @("Ardi", "unknown"),
    @("Ade", "unknown")
)

$row = 12
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

$ws.Range("B22").Select()
